$wb = $excel.ActiveWorkbook

# The workbook has a single worksheet, currently named "TEST_CASE".
$ws = $wb.Worksheets.Item("TEST_CASE")

# Correction of the worksheet name, to "TestCases"
$ws.Name = "TestCases"

# Reflect the last active selection on the sheet (was C3:D3, now D8)
$ws.Range("D8").Select()
